# Plasma_Gen_Q&A_list — 1/14 Femto 전달 문서
# Adds the two new H/W questions (rows 12 & 13 / sheet rows 15 & 16) to the
# "Plasma_Gen" sheet, shrinks row 9's height, and moves the active
# sheet/selection back to the "1월12일 meeting" tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Plasma_Gen sheet — fill in the two previously-empty rows (12 & 13)
# ---------------------------------------------------------------------
$gen = $wb.Worksheets.Item("Plasma_Gen")

# Row 15 (question #12)
$gen.Range("C15").Value = 43114
$gen.Range("D15").Value = "H/W"
$gen.Range("D15").HorizontalAlignment = -4108
$gen.Range("E15").Value = "OPEN"
$gen.Range("E15").HorizontalAlignment = -4108
$gen.Range("F15").Value = "External Power model에서 불필요한 기능 확인" + [char]10 + " - 전원 Key, Plasma Key, Volume Key" + [char]10 + " - Buzzer" + [char]10 + " - 충전 IC, Battery, Battery CON"
$gen.Range("F15").WrapText = $true

# Row 16 (question #13)
$gen.Range("C16").Value = 43114
$gen.Range("D16").Value = "H/W"
$gen.Range("D16").HorizontalAlignment = -4108
$gen.Range("E16").Value = "OPEN"
$gen.Range("E16").HorizontalAlignment = -4108
$gen.Range("F16").Value = "Battery model에서 RS-232 통신 지원 여부"

# Row heights — set *after* the content/wrap above so the explicit
# heights stick instead of being overridden by auto-fit.
$gen.Rows.Item(9).RowHeight = 49.5
$gen.Rows.Item(15).RowHeight = 66

# ---------------------------------------------------------------------
# 2. Sheet/selection bookkeeping — switch the active tab back to
#    "1월12일 meeting" and update each sheet's remembered selection.
# ---------------------------------------------------------------------
$meeting = $wb.Worksheets.Item("1월12일 meeting")
$lf = $wb.Worksheets.Item("Plasma_LF")

$lf.Activate()
$lf.Range("F19").Select()

$meeting.Activate()
$meeting.Range("I11").Select()
